# Add three new rows (117-119) to Sheet1, continuing the daily date series
# in column A and repeating the same B:J values found in the last existing
# row (116). Using Copy/Paste (instead of manually poking style properties)
# keeps the original cell style (s="2") intact instead of creating new,
# near-duplicate style entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 116
$newRowsCount = 3

$lastDate = $ws.Cells.Item($lastRow, 1).Value2
$srcRange = $ws.Range("A" + $lastRow + ":J" + $lastRow)

for ($i = 1; $i -le $newRowsCount; $i++) {
    $r = $lastRow + $i
    $dstRange = $ws.Range("A" + $r + ":J" + $r)

    # Copy the whole row (values + formatting) from the last row.
    $srcRange.Copy($dstRange)

    # Column A: date serial incremented by one day from the previous row.
    $ws.Cells.Item($r, 1).Value2 = $lastDate + $i
}
